$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (row 1): descriptive Spanish headers -> short codes ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case the Spanish prepositions (de/del/la/las/el/los/y) in place names ---
# e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga"
$ws.Range("B7").Value = 'Pabellón De Arteaga'
$ws.Range("B8").Value = 'Rincón De Romos'
$ws.Range("B9").Value = 'San Francisco De Los Romo'
$ws.Range("B35").Value = 'Amatenango De La Frontera'
$ws.Range("B38").Value = 'Bejucal De Ocampo'
$ws.Range("B40").Value = 'Benemérito De Las Américas'
$ws.Range("B48").Value = 'Chiapa De Corzo'
$ws.Range("B52").Value = 'Comitán De Domínguez'
$ws.Range("B78").Value = 'Marqués De Comillas'
$ws.Range("B79").Value = 'Mazapa De Madero'
$ws.Range("B85").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B94").Value = 'Salto De Agua'
$ws.Range("B95").Value = 'San Cristóbal De Las Casas'
$ws.Range("B128").Value = 'Coyame Del Sotol'
$ws.Range("B135").Value = 'Guadalupe Y Calvo'
$ws.Range("B138").Value = 'Hidalgo Del Parral'
$ws.Range("B150").Value = 'San Francisco Del Oro'
$ws.Range("B174").Value = 'San Juan De Sabinas'
$ws.Range("B188").Value = 'Villa De Álvarez'
$ws.Range("A190").Value = 'Ciudad De México'
$ws.Range("B194").Value = 'Cuajimalpa De Morelos'
$ws.Range("B221").Value = 'Nombre De Dios'
$ws.Range("B225").Value = 'Pánuco De Coronado'
$ws.Range("B230").Value = 'San Juan De Guadalupe'
$ws.Range("B231").Value = 'San Juan Del Río'
$ws.Range("B232").Value = 'San Luis Del Cordero'
$ws.Range("A242").Value = 'Estado De México'
$ws.Range("B242").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B245").Value = 'Almoloya De Alquisiras'
$ws.Range("B246").Value = 'Almoloya De Juárez'
$ws.Range("B251").Value = 'Atizapán De Zaragoza'
$ws.Range("B256").Value = 'Chapa De Mota'
$ws.Range("B261").Value = 'Coacalco De Berriozábal'
$ws.Range("B267").Value = 'Ecatepec De Morelos'
$ws.Range("B274").Value = 'Ixtapan De La Sal'
$ws.Range("B275").Value = 'Ixtapan Del Oro'
$ws.Range("B289").Value = 'Naucalpan De Juárez'
$ws.Range("B300").Value = 'San Antonio La Isla'
$ws.Range("B301").Value = 'San Felipe Del Progreso'
$ws.Range("B302").Value = 'San Martín De Las Pirámides'
$ws.Range("B304").Value = 'San Simón De Guerrero'
$ws.Range("B306").Value = 'Soyaniquilpan De Juárez'
$ws.Range("B316").Value = 'Tenango Del Valle'
$ws.Range("B327").Value = 'Tlalnepantla De Baz'
$ws.Range("B333").Value = 'Valle De Bravo'
$ws.Range("B334").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B335").Value = 'Villa De Allende'
$ws.Range("B336").Value = 'Villa Del Carbón'
$ws.Range("B349").Value = 'San Miguel De Allende'
$ws.Range("B350").Value = 'Apaseo El Alto'
$ws.Range("B351").Value = 'Apaseo El Grande'
$ws.Range("B359").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B363").Value = 'Jaral Del Progreso'
$ws.Range("B371").Value = 'Purísima Del Rincón'
$ws.Range("B375").Value = 'San Diego De La Unión'
$ws.Range("B377").Value = 'San Francisco Del Rincón'
$ws.Range("B379").Value = 'San Luis De La Paz'
$ws.Range("B381").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B383").Value = 'Silao De La Victoria'
$ws.Range("B388").Value = 'Valle De Santiago'
$ws.Range("B394").Value = 'Acapulco De Juárez'
$ws.Range("B397").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B398").Value = 'Alcozauca De Guerrero'
$ws.Range("B402").Value = 'Atenango Del Río'
$ws.Range("B403").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B405").Value = 'Atoyac De Álvarez'
$ws.Range("B406").Value = 'Ayutla De Los Libres'
$ws.Range("B409").Value = 'Buenavista De Cuéllar'
$ws.Range("B410").Value = 'Chilapa De Álvarez'
$ws.Range("B411").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B412").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B417").Value = 'Coyuca De Benítez'
$ws.Range("B418").Value = 'Coyuca De Catalán'
$ws.Range("B422").Value = 'Cuetzala Del Progreso'
$ws.Range("B423").Value = 'Cutzamala De Pinzón'
$ws.Range("B429").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B430").Value = 'Iguala De La Independencia'
$ws.Range("B432").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B435").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B438").Value = 'Mártir De Cuilapan'
$ws.Range("B451").Value = 'Taxco De Alarcón'
$ws.Range("B453").Value = 'Técpan De Galeana'
$ws.Range("B455").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B457").Value = 'Tixtla De Guerrero'
$ws.Range("B460").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B461").Value = 'Tlapa De Comonfort'
$ws.Range("B473").Value = 'Agua Blanca De Iturbide'
$ws.Range("B479").Value = 'Atotonilco De Tula'
$ws.Range("B480").Value = 'Atotonilco El Grande'
$ws.Range("B486").Value = 'Cuautepec De Hinojosa'
$ws.Range("B490").Value = 'Huasca De Ocampo'
$ws.Range("B494").Value = 'Huejutla De Reyes'
$ws.Range("B497").Value = 'Jacala De Ledezma'
$ws.Range("B503").Value = 'Mineral De La Reforma'
$ws.Range("B504").Value = 'Mineral Del Chico'
$ws.Range("B505").Value = 'Mineral Del Monte'
$ws.Range("B506").Value = 'Mixquiahuala De Juárez'
$ws.Range("B507").Value = 'Molango De Escamilla'
$ws.Range("B509").Value = 'Nopala De Villagrán'
$ws.Range("B510").Value = 'Omitlán De Juárez'
$ws.Range("B511").Value = 'Pachuca De Soto'
$ws.Range("B514").Value = 'Progreso De Obregón'
$ws.Range("B520").Value = 'Santiago De Anaya'
$ws.Range("B521").Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range("B525").Value = 'Tenango De Doria'
$ws.Range("B527").Value = 'Tepehuacán De Guerrero'
$ws.Range("B528").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B531").Value = 'Tezontepec De Aldama'
$ws.Range("B539").Value = 'Tula De Allende'
$ws.Range("B540").Value = 'Tulancingo De Bravo'
$ws.Range("B541").Value = 'Villa De Tezontepec'
$ws.Range("B544").Value = 'Zacualtipán De Ángeles'
$ws.Range("B545").Value = 'Zapotlán De Juárez'
$ws.Range("B550").Value = 'Ahualulco De Mercado'
$ws.Range("B557").Value = 'Atotonilco El Alto'
$ws.Range("B559").Value = 'Autlán De Navarro'
$ws.Range("B565").Value = 'Cañadas De Obregón'
$ws.Range("B572").Value = 'Concepción De Buenos Aires'
$ws.Range("B579").Value = 'Encarnación De Díaz'
$ws.Range("B586").Value = 'Huejuquilla El Alto'
$ws.Range("B587").Value = 'Ixtlahuacán Del Río'
$ws.Range("B591").Value = 'Jilotlán De Los Dolores'
$ws.Range("B596").Value = 'Lagos De Moreno'
$ws.Range("B605").Value = 'San Diego De Alejandría'
$ws.Range("B607").Value = 'San Juan De Los Lagos'
$ws.Range("B609").Value = 'San Martín De Bolaños'
$ws.Range("B610").Value = 'San Miguel El Alto'
$ws.Range("B611").Value = 'San Sebastián Del Oeste'
$ws.Range("B612").Value = 'Santa María De Los Ángeles'
$ws.Range("B614").Value = 'Talpa De Allende'
$ws.Range("B615").Value = 'Tamazula De Gordiano'
$ws.Range("B617").Value = 'Techaluta De Montenegro'
$ws.Range("B620").Value = 'Teocuitatlán De Corona'
$ws.Range("B621").Value = 'Tepatitlán De Morelos'
$ws.Range("B623").Value = 'Tizapán El Alto'
$ws.Range("B624").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B632").Value = 'Unión De San Antonio'
$ws.Range("B633").Value = 'Unión De Tula'
$ws.Range("B634").Value = 'Valle De Guadalupe'
$ws.Range("B635").Value = 'Valle De Juárez'
$ws.Range("B640").Value = 'Yahualica De González Gallo'
$ws.Range("B641").Value = 'Zacoalco De Torres'
$ws.Range("B644").Value = 'Zapotitlán De Vadillo'
$ws.Range("B645").Value = 'Zapotlán El Grande'
$ws.Range("B671").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B673").Value = 'Cojumatlán De Régules'
$ws.Range("B737").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B761").Value = 'Coatlán Del Río'
$ws.Range("B773").Value = 'Puente De Ixtla'
$ws.Range("B779").Value = 'Tetela Del Volcán'
$ws.Range("B780").Value = 'Tlaltizapán De Zapata'
$ws.Range("B786").Value = 'Zacualpan De Amilpas'
$ws.Range("B790").Value = 'Amatlán De Cañas'
$ws.Range("B794").Value = 'Ixtlán Del Río'
$ws.Range("B801").Value = 'Santa María Del Oro'
$ws.Range("B822").Value = 'Mier Y Noriega'
$ws.Range("B826").Value = 'San Nicolás De Los Garza'
$ws.Range("B831").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B836").Value = 'Ayoquezco De Aldama'
$ws.Range("B840").Value = 'Capulálpam De Méndez'
$ws.Range("B841").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B842").Value = 'Ciénega De Zimatlán'
$ws.Range("B845").Value = 'Coicoyán De Las Flores'
$ws.Range("B846").Value = 'Constancia Del Rosario'
$ws.Range("B848").Value = 'Cuilápam De Guerrero'
$ws.Range("B849").Value = 'Fresnillo De Trujano'
$ws.Range("B850").Value = 'Guadalupe De Ramírez'
$ws.Range("B852").Value = 'Guevea De Humboldt'
$ws.Range("B853").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B854").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B855").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B856").Value = 'Huautla De Jiménez'
$ws.Range("B858").Value = 'Ixtlán De Juárez'
$ws.Range("B859").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B869").Value = 'Mariscala De Juárez'
$ws.Range("B870").Value = 'Mártires De Tacubaya'
$ws.Range("B873").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B874").Value = 'Mixistlán De La Reforma'
$ws.Range("B878").Value = 'Nejapa De Madero'
$ws.Range("B879").Value = 'Oaxaca De Juárez'
$ws.Range("B880").Value = 'Ocotlán De Morelos'
$ws.Range("B881").Value = 'Pinotepa De Don Luis'
$ws.Range("B883").Value = 'Putla Villa De Guerrero'
$ws.Range("B884").Value = 'Reforma De Pineda'
$ws.Range("B886").Value = 'Rojas De Cuauhtémoc'
$ws.Range("B901").Value = 'San Antonino El Alto'
$ws.Range("B908").Value = 'San Baltazar Yatzachi El Bajo'
$ws.Range("B921").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B946").Value = 'San José Del Progreso'
$ws.Range("B955").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B964").Value = 'San Juan Del Río'
$ws.Range("B1016").Value = 'San Miguel Del Puerto'
$ws.Range("B1018").Value = 'San Miguel El Grande'
$ws.Range("B1037").Value = 'San Pablo Villa De Mitla'
$ws.Range("B1043").Value = 'San Pedro El Alto'
$ws.Range("B1058").Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range("B1059").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B1073").Value = 'Santa Ana Del Valle'
$ws.Range("B1083").Value = 'Santa Cruz De Bravo'
$ws.Range("B1088").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B1093").Value = 'Santa Inés De Zaragoza'
$ws.Range("B1094").Value = 'Santa Inés Del Monte'
$ws.Range("B1096").Value = 'Santa Lucía Del Camino'
$ws.Range("B1106").Value = 'Santa María Del Rosario'
$ws.Range("B1107").Value = 'Santa María Del Tule'
$ws.Range("B1114").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B1143").Value = 'Santiago Del Río'
$ws.Range("B1175").Value = 'Santo Domingo De Morelos'
$ws.Range("B1191").Value = 'Sitio De Xitlapehua'
$ws.Range("B1193").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B1194").Value = 'Tanetze De Zaragoza'
$ws.Range("B1196").Value = 'Tataltepec De Valdés'
$ws.Range("B1197").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B1198").Value = 'Teotitlán De Flores Magón'
$ws.Range("B1199").Value = 'Teotitlán Del Valle'
$ws.Range("B1202").Value = 'Tlacolula De Matamoros'
$ws.Range("B1204").Value = 'Tlalixtac De Cabrera'
$ws.Range("B1205").Value = 'Totontepec Villa De Morelos'
$ws.Range("B1208").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B1209").Value = 'Villa De Etla'
$ws.Range("B1210").Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range("B1211").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B1212").Value = 'Villa De Zaachila'
$ws.Range("B1215").Value = 'Villa Sola De Vega'
$ws.Range("B1216").Value = 'Villa Talea De Castro'
$ws.Range("B1222").Value = 'Zimatlán De Álvarez'
$ws.Range("B1243").Value = 'Ayotoxco De Guerrero'
$ws.Range("B1246").Value = 'Chalchicomula De Sesma'
$ws.Range("B1256").Value = 'Chila De La Sal'
$ws.Range("B1265").Value = 'Cuapiaxtla De Madero'
$ws.Range("B1268").Value = 'Cuayuca De Andrade'
$ws.Range("B1281").Value = 'Huehuetlán El Chico'
$ws.Range("B1282").Value = 'Huehuetlán El Grande'
$ws.Range("B1285").Value = 'Huitzilan De Serdán'
$ws.Range("B1287").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B1290").Value = 'Izúcar De Matamoros'
$ws.Range("B1298").Value = 'Los Reyes De Juárez'
$ws.Range("B1308").Value = 'Palmar De Bravo'
$ws.Range("B1317").Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Range("B1328").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B1332").Value = 'San Salvador El Seco'
$ws.Range("B1333").Value = 'San Salvador El Verde'
$ws.Range("B1338").Value = 'Tecali De Herrera'
$ws.Range("B1346").Value = 'Tepanco De López'
$ws.Range("B1347").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B1353").Value = 'Tepexi De Rodríguez'
$ws.Range("B1355").Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range("B1356").Value = 'Tetela De Ocampo'
$ws.Range("B1357").Value = 'Teteles De Avila Castillo'
$ws.Range("B1362").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B1377").Value = 'Xayacatlán De Bravo'
$ws.Range("B1394").Value = 'Amealco De Bonfil'
$ws.Range("B1396").Value = 'Cadereyta De Montes'
$ws.Range("B1402").Value = 'Jalpan De Serra'
$ws.Range("B1403").Value = 'Landa De Matamoros'
$ws.Range("B1406").Value = 'Pinal De Amoles'
$ws.Range("B1409").Value = 'San Juan Del Río'
$ws.Range("B1420").Value = 'Armadillo De Los Infante'
$ws.Range("B1421").Value = 'Axtla De Terrazas'
$ws.Range("B1427").Value = 'Ciudad Del Maíz'
$ws.Range("B1437").Value = 'Mexquitic De Carmona'
$ws.Range("B1443").Value = 'San Ciro De Acosta'
$ws.Range("B1448").Value = 'Santa María Del Río'
$ws.Range("B1450").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B1458").Value = 'Tanquián De Escobedo'
$ws.Range("B1462").Value = 'Villa De Arista'
$ws.Range("B1463").Value = 'Villa De Arriaga'
$ws.Range("B1464").Value = 'Villa De Guadalupe'
$ws.Range("B1465").Value = 'Villa De La Paz'
$ws.Range("B1466").Value = 'Villa De Ramos'
$ws.Range("B1467").Value = 'Villa De Reyes'
$ws.Range("B1502").Value = 'Nacozari De García'
$ws.Range("B1517").Value = 'Jalpa De Méndez'
$ws.Range("B1557").Value = 'Soto La Marina'
$ws.Range("B1571").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1574").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1577").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1580").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1581").Value = 'San Pablo Del Monte'
$ws.Range("B1584").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1587").Value = 'Tetla De La Solidaridad'
$ws.Range("B1608").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1612").Value = 'Amatlán De Los Reyes'
$ws.Range("B1624").Value = 'Boca Del Río'
$ws.Range("B1626").Value = 'Camarón De Tejeda'
$ws.Range("B1630").Value = 'Castillo De Teayo'
$ws.Range("B1632").Value = 'Cazones De Herrera'
$ws.Range("B1650").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1651").Value = 'Cosautlán De Carvajal'
$ws.Range("B1667").Value = 'Hueyapan De Ocampo'
$ws.Range("B1668").Value = 'Ignacio De La Llave'
$ws.Range("B1672").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B1673").Value = 'Ixhuatlán De Madero'
$ws.Range("B1674").Value = 'Ixhuatlán Del Café'
$ws.Range("B1675").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B1686").Value = 'Juchique De Ferrer'
$ws.Range("B1689").Value = 'Landero Y Coss'
$ws.Range("B1692").Value = 'Lerdo De Tejada'
$ws.Range("B1696").Value = 'Martínez De La Torre'
$ws.Range("B1698").Value = 'Medellín De Bravo'
$ws.Range("B1702").Value = 'Mixtla De Altamirano'
$ws.Range("B1704").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1712").Value = 'Ozuluama De Mascareñas'
$ws.Range("B1716").Value = 'Paso De Ovejas'
$ws.Range("B1717").Value = 'Paso Del Macho'
$ws.Range("B1721").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1729").Value = 'Sayula De Alemán'
$ws.Range("B1733").Value = 'Soledad De Doblado'
$ws.Range("B1758").Value = 'Tlacotepec De Mejía'
$ws.Range("B1770").Value = 'Vega De Alatorre'
$ws.Range("B1781").Value = 'Zozocolco De Hidalgo'
$ws.Range("B1798").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B1800").Value = 'Concepción Del Oro'
$ws.Range("B1811").Value = 'Jiménez Del Teul'
$ws.Range("B1818").Value = 'Mezquital Del Oro'
$ws.Range("B1822").Value = 'Moyahua De Estrada'
$ws.Range("B1823").Value = 'Nochistlán De Mejía'
$ws.Range("B1824").Value = 'Noria De Ángeles'
$ws.Range("B1835").Value = 'Teúl De González Ortega'
$ws.Range("B1836").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1839").Value = 'Villa De Cos'

# --- Floating point percentage recalculation artifacts (1-ULP bumps) ---
$ws.Range("D353").Value = 0.009593094944512949
$ws.Range("D527").Value = 0.0009617755856966709
$ws.Range("D535").Value = 0.0009617755856966709
$ws.Range("D855").Value = 0.0009617755856966709
$ws.Range("D1226").Value = 0.0009617755856966709
$ws.Range("D1238").Value = 0.0009617755856966709
$ws.Range("D1404").Value = 0.0009617755856966709
$ws.Range("D1416").Value = 0.0009617755856966709
$ws.Range("D1437").Value = 0.0009617755856966709
$ws.Range("D1687").Value = 0.0009617755856966709
$ws.Range("D1721").Value = 0.0009617755856966709

# --- Remove trailing footnote rows (1848:1852) and shrink used range to A1:D1846 ---
$ws.Range("A1848:A1852").EntireRow.Delete()